$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / extend the data rows ---
# Columns: A = sr (number), B = gitlab_project_namespace, C = project_to_import, D = github_username
# Write values in the order that reproduces the author's original shared-string table order:
# repo-migration, code-migration, casa-build-utils, app-n-pak, almatasks, casa6, casashell
$ws.Range("B2").Value = "repo-migration"
$ws.Range("D2").Value = "code-migration"
$ws.Range("C4").Value = "casa-build-utils"
$ws.Range("C3").Value = "app-n-pak"
$ws.Range("C2").Value = "almatasks"
$ws.Range("C5").Value = "casa6"
$ws.Range("C6").Value = "casashell"

# Remaining repeated namespace / username values
$ws.Range("B3").Value = "repo-migration"
$ws.Range("D3").Value = "code-migration"

# New rows 4-6
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "repo-migration"
$ws.Range("D4").Value = "code-migration"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "repo-migration"
$ws.Range("D5").Value = "code-migration"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "repo-migration"
$ws.Range("D6").Value = "code-migration"

# --- Formatting: project_to_import / github_username columns get a (no-fill) style tag ---
# applied to every data cell in C/D except C5, matching the source workbook.
$ws.Range("C2:D4").Interior.ColorIndex = -4142
$ws.Range("D5").Interior.ColorIndex = -4142
$ws.Range("C6:D6").Interior.ColorIndex = -4142

# --- Selection state ---
$ws.Range("B6").Select()
